# Update countries & provincias Spain
# - Refresh the "last updated" timestamp
# - Apply the latest daily COVID-19 figures for each country row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 2 de Abril de 2020 a las 17:50"

$ws.Range("B4").Value = 226378
$ws.Range("C4").Value = 11375
$ws.Range("D4").Value = 10265
$ws.Range("E4").Value = 210779
$ws.Range("F4").Value = 5403
$ws.Range("G4").Value = 232
$ws.Range("H4").Value = 5334

$ws.Range("B12").Value = 18475
$ws.Range("C12").Value = 707
$ws.Range("E12").Value = 13940
$ws.Range("G12").Value = 34
$ws.Range("H12").Value = 522

$ws.Range("B16").Value = 11027
$ws.Range("C16").Value = 316
$ws.Range("E16").Value = 9120

$ws.Range("D17").Value = 1891
$ws.Range("E17").Value = 8110
$ws.Range("G17").Value = 17
$ws.Range("H17").Value = 131

$ws.Range("B24").Value = 5124
$ws.Range("C24").Value = 247
$ws.Range("E24").Value = 5044
$ws.Range("G24").Value = 4
$ws.Range("H24").Value = 48

$ws.Range("B25").Value = 3805
$ws.Range("C25").Value = 216
$ws.Range("D25").Value = 67
$ws.Range("E25").Value = 3694
$ws.Range("G25").Value = 5
$ws.Range("H25").Value = 44

$ws.Range("F27").Value = 126

$ws.Range("B35").Value = 2487
$ws.Range("C35").Value = 168
$ws.Range("D35").Value = 80
$ws.Range("E35").Value = 2377
$ws.Range("F35").Value = 31
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 30

$ws.Range("B36").Value = 2386
$ws.Range("C36").Value = 268
$ws.Range("D36").Value = 107
$ws.Range("E36").Value = 2247
$ws.Range("F36").Value = 12
$ws.Range("G36").Value = 5
$ws.Range("H36").Value = 32

$ws.Range("B37").Value = 2384
$ws.Range("D37").Value = 472
$ws.Range("E37").Value = 1855
$ws.Range("F37").Value = 69
$ws.Range("H37").Value = 57

$ws.Range("B38").Value = 2341
$ws.Range("C38").Value = 343
$ws.Range("D38").Value = 177
$ws.Range("E38").Value = 2096
$ws.Range("G38").Value = 10
$ws.Range("H38").Value = 68

$ws.Range("B43").Value = 1514
$ws.Range("C43").Value = 99
$ws.Range("D43").Value = 61
$ws.Range("E43").Value = 1400
$ws.Range("F43").Value = 91
$ws.Range("G43").Value = 2
$ws.Range("H43").Value = 53

$ws.Range("E71").Value = 482
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 16

$ws.Range("B76").Value = 455
$ws.Range("C76").Value = 32
$ws.Range("D76").Value = 5
$ws.Range("E76").Value = 436
$ws.Range("F76").Value = 10
$ws.Range("G76").Value = 2
$ws.Range("H76").Value = 14

$ws.Range("B77").Value = 428
$ws.Range("C77").Value = 38
$ws.Range("D77").Value = 10
$ws.Range("E77").Value = 403
$ws.Range("F77").Value = 12
$ws.Range("G77").Value = 1
$ws.Range("H77").Value = 15

$ws.Range("B78").Value = 426
$ws.Range("C78").Value = 26
$ws.Range("E78").Value = 420
$ws.Range("F78").Value = 3
$ws.Range("H78").Value = 1

$ws.Range("E83").Value = 318
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 10

$ws.Range("F88").Value = 11

$ws.Range("B113").Value = 131
$ws.Range("C113").Value = 14
$ws.Range("E113").Value = 105

$ws.Range("E115").Value = 114
$ws.Range("G115").Value = 1
$ws.Range("H115").Value = 8

$ws.Range("C153").Value = 4
$ws.Range("D153").Value = 0
$ws.Range("E153").Value = 19

$ws.Range("B154").Value = 20
$ws.Range("D154").Value = 2
$ws.Range("E154").Value = 17
$ws.Range("G154").Value = 0
$ws.Range("H154").Value = 1

$ws.Range("D155").Value = 0
$ws.Range("E155").Value = 15
$ws.Range("G155").Value = 1
$ws.Range("H155").Value = 4

$ws.Range("B156").Value = 19
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 13
$ws.Range("E156").Value = 6

$ws.Range("C157").Value = 3
$ws.Range("E157").Value = 18
$ws.Range("H157").Value = 0

$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 0
$ws.Range("H158").Value = 1

$ws.Range("B159").Value = 18
$ws.Range("C159").Value = 2
$ws.Range("D159").Value = 1

$ws.Range("B160").Value = 17
$ws.Range("D160").Value = 0
$ws.Range("E160").Value = 17

$ws.Range("D161").Value = 1
$ws.Range("H161").Value = 0
